$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 140 (pushes "Participants" etc. down by one row),
#    inheriting the yellow-highlight style (style 5) from row 139 above it.
$ws.Range("A140").EntireRow.Insert()
$ws.Range("A140").Value = "CourtOfficerNotes"
$ws.Range("B140").Value = "cares\Courts.xlsx"
$ws.Range("C140").Value = "CourtOfficerNotes"
$ws.Range("D140").Value = 1

# 1b) Fill in the previously-blank "iteration count" column for the four
#     Court rows (CourtCase/CourtHearings/CourtWarrents/CourtWorkItems),
#     now sitting at rows 144:147 after the insert above.
$ws.Range("D144:D147").Value = 1

# 2) After the insert above, "CourtWorkItems" is now row 147.
#    Insert 3 new rows right after it (new rows 148:150) for the
#    DocumentDistributions / CourtReports / ReleasesOfInformationAndConsentForm entries.
$ws.Range("A148:A150").EntireRow.Insert()

$ws.Range("A148").Value = "DocumentDistributions"
$ws.Range("B148").Value = "cares\Courts.xlsx"
$ws.Range("C148").Value = "DocumentDistributions"
$ws.Range("D148").Value = 1

$ws.Range("A149").Value = "CourtReports"
$ws.Range("B149").Value = "cares\Courts.xlsx"
$ws.Range("C149").Value = "CourtReports"
$ws.Range("D149").Value = 1

$ws.Range("A150").Value = "ReleasesOfInformationAndConsentForm"
$ws.Range("B150").Value = "cares\Courts.xlsx"
$ws.Range("C150").Value = "ReleasesOfInformationAndConsentForm"
$ws.Range("D150").Value = 1

# 3) Append 4 new rows (175:178) after the current last row (174, "SubmitToCalSAWS")
#    for the new Placement-preference entries.
$ws.Range("A175").Value = "TribalPlacementPreferences"
$ws.Range("B175").Value = "cares\Placement.xlsx"
$ws.Range("C175").Value = "TribalPlacementPreferences"
$ws.Range("D175").Value = 1

$ws.Range("A176").Value = "ParentGuardPlacementPreferences"
$ws.Range("B176").Value = "cares\Placement.xlsx"
$ws.Range("C176").Value = "ParentGuardPlacementPreferences"
$ws.Range("D176").Value = 1

$ws.Range("A177").Value = "YouthPlacementPreferences"
$ws.Range("B177").Value = "cares\Placement.xlsx"
$ws.Range("C177").Value = "YouthPlacementPreferences"
$ws.Range("D177").Value = 1

$ws.Range("A178").Value = "PlacementNeeds"
$ws.Range("B178").Value = "cares\Placement.xlsx"
$ws.Range("C178").Value = "PlacementNeeds"
$ws.Range("D178").Value = 1

# Match the saved cursor/selection position recorded in the target file.
$ws.Application.ActiveWindow.ScrollRow = 136
$ws.Range("G147").Select()
